$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Locatie op schema" column (G) entirely - it only had 3 sparse
# values (G1 header, G2 "1:ESP32", G11 "4:Power") that are no longer used.
$ws.Columns("G").Delete()

# New components added to the BOM: two Molex connectors, each get their own
# row with the usual Aantal/Leverancier/Link columns filled in.
$ws.Range("B13").Value = "5 pin molex connector"
$ws.Range("C13").Value = 2
$ws.Range("D13").Value = "mouser"
$ws.Range("E13").Value = "link"
$ws.Hyperlinks.Add($ws.Range("E13"), "https://www.mouser.be/ProductDetail/5-pin-molex-connector")
$ws.Range("E13").Style = "Hyperlink"

$ws.Range("B14").Value = "2 pin molex connector"
$ws.Range("C14").Value = 1
$ws.Range("D14").Value = "mouser"
$ws.Range("E14").Value = "link"
$ws.Hyperlinks.Add($ws.Range("E14"), "https://www.mouser.be/ProductDetail/2-pin-molex-connector")
$ws.Range("E14").Style = "Hyperlink"

# C7 row now also covers C4 and C5 -> update the reference designator text.
$ws.Range("A7").Value = "C4,C5,C7"

# Leave the selection where the author left it before saving.
$ws.Range("A15").Select()
